$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "244.65"
Set-TextValue "E2" "-0.94%"
Set-TextValue "D3" "27.10"
Set-TextValue "E3" "2.56%"
Set-TextValue "D4" "5.068"
Set-TextValue "E4" "-0.30%"
Set-TextValue "D5" "0.05688"
Set-TextValue "E5" "1.29%"
Set-TextValue "D7" "0.8203"
Set-TextValue "E7" "0.89%"
Set-TextValue "D8" "0.8377"
Set-TextValue "E8" "-0.93%"
Set-TextValue "D9" "0.1329"
Set-TextValue "E9" "-1.30%"
Set-TextValue "D11" "0.02858"
Set-TextValue "E11" "-0.77%"
Set-TextValue "D12" "0.09400"
Set-TextValue "E12" "-0.11%"
Set-TextValue "D13" "0.001519"
Set-TextValue "E13" "-0.14%"
Set-TextValue "D14" "0.04114"
Set-TextValue "E14" "-11.74%"
Set-TextValue "D15" "0.0005978"
Set-TextValue "E15" "-0.38%"
Set-TextValue "D16" "0.006136"
Set-TextValue "E16" "-0.82%"
Set-TextValue "D17" "3.511"
Set-TextValue "E17" "-2.17%"
Set-TextValue "D18" "3.001"
Set-TextValue "E18" "-0.52%"
Set-TextValue "D19" "2.226"
Set-TextValue "E19" "5.08%"
Set-TextValue "D20" "0.3150"
Set-TextValue "E20" "-0.22%"
Set-TextValue "D21" "0.03190"
Set-TextValue "E21" "-0.48%"
Set-TextValue "D22" "0.1295"
Set-TextValue "E22" "-1.90%"
Set-TextValue "D23" "3.574"
Set-TextValue "E23" "-4.68%"
Set-TextValue "E24" "1.75%"
Set-TextValue "D25" "0.001218"
Set-TextValue "E25" "-2.42%"
Set-TextValue "D26" "0.003951"
Set-TextValue "E26" "-14.12%"
Set-TextValue "D27" "0.00009799"
Set-TextValue "E27" "2.09%"
Set-TextValue "E28" "-0.06%"
Set-TextValue "D40" "0.03695"
Set-TextValue "E40" "0.42%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1055"
Set-TextValue "E41" "-0.43%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.003427"
Set-TextValue "E42" "-44.08%"
Set-TextValue "D43" "0.002350"
Set-TextValue "E43" "-6.00%"
Set-TextValue "D44" "0.009372"
Set-TextValue "E44" "5.29%"
Set-TextValue "D45" "0.00005198"
Set-TextValue "E45" "-1.85%"
Set-TextValue "E47" "-32.33%"
Set-TextValue "E48" "2.42%"
